$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "['MEC-3B-Fresagem', -, -, -]"

# Row 4
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "['MEC-3B-Fresagem', -, -, -]"

# Row 6
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "['MEC-3B-Fresagem', -, -, -]"

# Row 7
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "['MEC-3B-Fresagem', -, -, -]"

# Row 10
$ws.Range("F10").Value = "['MEC-2A-Tornearia', -, -, -]"

# Row 11
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "[-, 'MEC-3A-Fresagem', -, -]"

# Row 12
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "[-, 'MEC-3A-Fresagem', -, -]"

# Row 14
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "[-, 'MEC-3A-Fresagem', 'MEC-2A-Tornearia', -]"

# Row 15
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "[-, 'MEC-3A-Fresagem', 'MEC-2A-Tornearia', -]"

# Row 16
$ws.Range("F16").Value = "[-, -, 'MEC-2A-Tornearia', -]"

# Row 18
$ws.Range("B18").Value = "['MEC-2NA-Fresagem', -, -, 'MEC-2NA-CAD/CAM']"
$ws.Range("C18").Value = "-"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "[-, -, 'MEC-2NB-Fresagem', -]"

# Row 19
$ws.Range("B19").Value = "['MEC-2NA-Fresagem', -, -, 'MEC-2NA-CAD/CAM']"
$ws.Range("C19").Value = "-"
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "[-, -, 'MEC-2NB-Fresagem', -]"

# Row 20
$ws.Range("B20").Value = "['MEC-2NA-Fresagem', -, -, 'MEC-2NA-CAD/CAM']"
$ws.Range("C20").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "[-, -, 'MEC-2NB-Fresagem', -]"

# Row 21
$ws.Range("B21").Value = "['MEC-2NA-Fresagem', -, -, 'MEC-2NA-CAD/CAM']"
$ws.Range("C21").Value = "-"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "[-, -, 'MEC-2NB-Fresagem', -]"
